$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "Amr Al Dhaheri"
$summary.Range("B4").Value = 4723.3
$summary.Range("B6").Value = 291050
$summary.Range("B7").Value = 219183
$summary.Range("B8").Value = 71867
$summary.Range("B9").Value = 1.33

# ---------------------------------------------------------------
# Sheet: Assets
# ---------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")
$assets.Range("B2").Value = "Premium Car"
$assets.Range("C2").Value = 283151
$assets.Range("C3").Value = 7899
$assets.Range("C4").Value = 291050

# ---------------------------------------------------------------
# Sheet: Liabilities
# ---------------------------------------------------------------
$liab = $wb.Worksheets.Item("Liabilities")

# Update existing Auto Loans row (row 2)
$liab.Range("C2").Value = 169891
$liab.Range("D2").Value = 3539
$liab.Range("E2").Value = 4

# Insert a new row before the old row 3 (Credit Cards), pushing it (and the
# TOTAL LIABILITIES row) down by one, so a "Personal Loans" row can be
# inserted at row 3.
$liab.Rows.Item(3).Insert()

$liab.Range("A3").Value = "Personal Loans"
$liab.Range("B3").Value = "Personal Loan"
$liab.Range("C3").Value = 26742
$liab.Range("D3").Value = 557
$liab.Range("E3").Value = 4

# Copy the formatting used by the other data rows (e.g. row 2) onto the
# newly inserted row so it matches the sheet's existing style.
$liab.Range("A2:E2").Copy()
$liab.Range("A3:E3").PasteSpecial(-4122)  # xlPasteFormats

# Now update what used to be row 3 (Credit Cards), now shifted to row 4.
$liab.Range("C4").Value = 22550
$liab.Range("D4").Value = 1128
$liab.Range("E4").Value = 1

# And the TOTAL LIABILITIES row, now shifted to row 5.
$liab.Range("C5").Value = 219183
